$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8589433431625366
$ws.Range("B1").Value = 2.970970153808594
$ws.Range("C1").Value = 4.584014892578125
$ws.Range("D1").Value = 2.848881483078003
$ws.Range("E1").Value = 1.428370833396912
